$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update master data values as part of the 2nd May data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update the sheet view selection: select from row 34 through the end of the
# sheet (this also clears the previous topLeftCell scroll position).
$ws.Rows("34:1048576").Select()
